# Split "natural gas nonpeaker" into "natural gas steam turbine" and
# "natural gas combined cycle", and remove the "(Boolean)" quality-tier
# suffix from the value-column header, replacing it with a dedicated
# units sub-header in the BDSBaPCF sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BDSBaPCF")

# Make room for the new "natural gas combined cycle" row right after the
# (renamed) "natural gas nonpeaker" row -> "natural gas steam turbine".
# Inserting a whole row shifts every row below it down by one and keeps
# all relative formulas (=B2, =B6, =B11, ...) pointing at the correct,
# shifted cells automatically.
$ws.Rows.Item(4).Insert()

# Row 3: rename "natural gas nonpeaker" -> "natural gas steam turbine"
$ws.Range("A3").Value = "natural gas steam turbine"

# Row 4 (new): "natural gas combined cycle", also bids at peak (1)
$ws.Range("A4").Value = "natural gas combined cycle"
$ws.Range("B4").Value = 1

# Header row: drop the "(Boolean)" qualifier from the value header and
# add a new italic units sub-header in column A.
$ws.Range("B1").Value = "Do Suppliers Bid at Peak Capacity Factors"
$ws.Range("A1").Value = "Unit: dimensionless (Boolean)"
$ws.Range("A1").Font.Italic = $true

# The shorter header text now wraps onto fewer lines, so the row is shorter.
$ws.Rows.Item(1).RowHeight = 30

# Column width tweaks to fit the new header text.
$ws.Columns.Item(1).ColumnWidth = 30.140625
$ws.Columns.Item(2).ColumnWidth = 23.140625

# Explicit portrait page orientation (now specified on this sheet).
$ws.PageSetup.Orientation = 1
